$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2 and D3 hold "UserName" values for the two sample rows. The column is
# formatted as Text ("@"), but the stored cell values are actual numbers
# (t="n"). Assigning .Value/.Value2 directly on a text-formatted cell makes
# Excel store the input as a text string instead, so the number format is
# temporarily switched to a numeric format, the new numeric value is written,
# and then the original Text format is restored so the cell keeps its
# original look (style) while the underlying value stays numeric.
$ws.Range("D2:D3").NumberFormat = "0"

$ws.Range("D2").Value2 = 5697768474
$ws.Range("D3").Value2 = 8857887012

$ws.Range("D2:D3").NumberFormat = "@"
